$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to retain text representation so that
# values such as "1.005" are not auto-coerced into numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.211.68'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.653.69'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = '218.18'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('D6').Value = '0.5196'
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '0.2647'
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('D9').Value = '0.06302'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').Value = '0.07735'
$ws.Range('D12').Value = '1.657.23'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = '4.420'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('D14').Value = '0.5447'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '0.0₅8194'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = '64.65'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').Value = '26.214.56'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D19').Value = '4.678'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').Value = '190.81'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D22').Value = '6.182'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').Value = '138.89'
$ws.Range('E24').Value = '  -3.06%  '
$ws.Range('D25').Value = '0.1242'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').Value = '7.277'
$ws.Range('E26').Value = '  -1.85%  '
$ws.Range('D27').Value = '16.07'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').Value = '1.416'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '0.06065'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').Value = '3.544'
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('D32').Value = '3.355'
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').Value = '0.9839'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').Value = '2.414'
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').Value = '2.770'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = '0.5927'
$ws.Range('E37').Value = '  +4.70%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').Value = '5.957'
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('D40').Value = '0.8622'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').Value = '1.057.84'
$ws.Range('E41').Value = '  +2.43%  '
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').Value = '99.77'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '1.794.78'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '0.0₈109'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').Value = '57.24'
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').Value = '8.055'
$ws.Range('E48').Value = '  -0.69%  '
$ws.Range('D49').Value = '0.05179'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').Value = '1.469'
$ws.Range('E50').Value = '  +5.08%  '
$ws.Range('D51').Value = '0.4233'
$ws.Range('E51').Value = '  +0.49%  '

# Restore the original (default) style now that the text values are set.
$dataRange.Style = $origStyle
